$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "Case" query stored in B2 had a stray leading space before the final
# "order By ..." line, which made the limit-10 reading function only match
# the first tab's query string and skip the others. Remove that leading
# space so all tabs are read consistently.
$caseQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
WHERE s.study_acronym IN ["C"]  
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
order By ss.study_subject_id ASC LIMIT 100 
'@

$ws.Range("B2").Value = $caseQuery

# Update the active selection left behind by the editing session.
$ws.Range("C2").Select()
